# Clean up the "length group" labels in column A: strip the dotted leaders
# that were used as visual fillers in the original OCR/typed table, leaving
# plain "NN-NN" (or "NN- NN ") text labels. Numeric data (columns B and C)
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list (row, label) -- order matters because it controls the order
# in which newly-introduced label strings land in the shared-string table.
$labelOrder = @(
    @(1,  "Length group in feet"),
    @(2,  "0- 10"),
    @(4,  "16- 20 "),
    @(6,  "26- 30"),
    @(7,  "31- 35 "),
    @(8,  "36- 40 "),
    @(9,  "41- 45 "),
    @(10, "46- 50 "),
    @(11, "51- 55 "),
    @(12, "56- 60 "),
    @(13, "61- 65 "),
    @(14, "66- 70"),
    @(15, "71- 75 "),
    @(16, "76- 80 "),
    @(19, "91- 95 "),
    @(21, "101-105"),
    @(22, "106-110"),
    @(23, "111-115"),
    @(25, "121-125"),
    @(26, "126-130 "),
    @(27, "131-135"),
    @(28, "136-140 "),
    @(29, "141-145 "),
    @(30, "146-150 "),
    @(31, "151-155"),
    @(32, "156-160"),
    @(33, "161-165 "),
    @(34, "166-170 "),
    @(35, "171-175"),
    @(3,  "11-15")
)

foreach ($pair in $labelOrder) {
    $row = $pair[0]
    $text = $pair[1]
    $ws.Cells.Item($row, 1).Value = $text
}

# Labels that were already "clean" (no dotted leader) keep their value but
# still need to be re-pointed at the (reordered) shared-string table; simply
# re-assigning the same text is enough to normalize them.
$unchangedLabels = @{
    5  = "21-25"
    17 = "81- 85"
    18 = "86-90"
    20 = "96-100"
    24 = "116-120"
    36 = "176-180"
    37 = "181 and over"
    38 = "Total"
    39 = "Total check"
}
foreach ($row in $unchangedLabels.Keys) {
    $ws.Cells.Item($row, 1).Value = $unchangedLabels[$row]
}

# The "11-15" label row now stores its text value with an explicit Text
# number format (it previously looked numeric-ish with the dotted leader).
$ws.Range("A3").NumberFormat = "@"

# Update the saved selection to match the author's last position.
$ws.Range("D34").Select()
